$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and 1h volume-change (E) values for the refreshed rows
$ws.Range("D2").Value = "26.896.45"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "1.553.33"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.72"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0859"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "1.774.51"
$ws.Range("E12").Value = "  +1.49%  "
$ws.Range("D13").Value = "1.541.68"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("E15").Value = "  +2.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "26.885.09"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.70%  "
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.62%  "
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("E26").Value = "  +2.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0464"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "1.418.57"
$ws.Range("E33").Value = "  +4.81%  "
$ws.Range("E34").Value = "  +3.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.962"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.05%  "
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.522"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.807"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.64%  "
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("E44").Value = "  +3.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.34%  "
$ws.Range("E46").Value = "  +1.53%  "
$ws.Range("D47").Value = "1.688.48"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0519"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.01%  "

# Row 50/51 restructure: the USDD row is removed, Algorand moves down to row 51
# with refreshed data, and a new BabyDogeCoin row is inserted at row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0101"
$ws.Range("E50").Value = "  +11.54%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0961"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.73%  "

